$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chapter 12 rows -------------------------------------------------------
$ws.Range("A18").Value = 12
$ws.Range("B18").Value = 2159
$ws.Range("C18").Value = "MVAfactcarm"
$ws.Range("D18").Value = 2159
$ws.Range("E18").Value = "Ready"
$ws.Range("F18").Value = 1398
$ws.Range("G18").Value = "X"
$ws.Range("H18").Value = "-"
$ws.Range("I18").Value = "Ready"
$ws.Range("K18").Value = "ok"

$ws.Range("A19").Value = 12
$ws.Range("B19").Value = 1211
$ws.Range("C19").Value = "MVAfacthous"
$ws.Range("D19").Value = 1211
$ws.Range("E19").Value = "Ready"
$ws.Range("F19").Value = 1672
$ws.Range("G19").Value = "X"
$ws.Range("H19").Value = "-"
$ws.Range("I19").Value = "-"
$ws.Range("K19").Value = "doesn't match | graph 4: mirrored image"

# Chapter 13 rows -------------------------------------------------------
$ws.Range("A20").Value = 13
$ws.Range("B20").Value = 1201
$ws.Range("C20").Value = "MVAclus8p"
$ws.Range("D20").Value = 1201
$ws.Range("E20").Value = "Ready"
$ws.Range("F20").Value = 1664
$ws.Range("G20").Value = "X"
$ws.Range("H20").Value = "-"
$ws.Range("I20").Value = "Ready"
$ws.Range("K20").Value = "slight differences | first graph"

$ws.Range("A21").Value = 13
$ws.Range("B21").Value = 1202
$ws.Range("C21").Value = "MVAclusbank"
$ws.Range("D21").Value = 1202
$ws.Range("E21").Value = "Ready"
$ws.Range("F21").Value = 1665
$ws.Range("G21").Value = "X"
$ws.Range("H21").Value = "-"
$ws.Range("I21").Value = "Ready"
$ws.Range("K21").Value = "graphs don't match"

$ws.Range("A22").Value = 13
$ws.Range("B22").Value = 1205
$ws.Range("C22").Value = "MVAclusfood"
$ws.Range("D22").Value = 1205
$ws.Range("E22").Value = "Ready"
$ws.Range("F22").Value = 1667
$ws.Range("G22").Value = "X"
$ws.Range("H22").Value = "-"
$ws.Range("I22").Value = "Ready"
$ws.Range("K22").Value = "ok"

$ws.Range("A23").Value = 13
$ws.Range("B23").Value = 1204
$ws.Range("C23").Value = "MVAclusbh"
$ws.Range("D23").Value = 1204
$ws.Range("E23").Value = "Ready"
$ws.Range("F23").Value = 1666
$ws.Range("G23").Value = "X"
$ws.Range("H23").Value = "-"
$ws.Range("I23").Value = "Ready"
$ws.Range("K23").Value = "ok"

# Chapter 14 starts with just the chapter number -------------------------
$ws.Range("A24").Value = 14

# The "Review Date" column (J) already holds the text "08.09.2016" a few
# rows up (J17); copy it down instead of typing it so Excel keeps storing
# it as plain text (typing it fresh gets auto-parsed into a date serial).
for ($r = 18; $r -le 23; $r++) {
    $ws.Range("J17").Copy() | Out-Null
    $ws.Range("J$r").PasteSpecial(-4104) | Out-Null
}

# K19 carries the "highlighted issue" formatting used elsewhere in the
# sheet for rows where the graphs don't match (e.g. K6/K10/K12/K13) -----
$ws.Range("K6").Copy() | Out-Null
$ws.Range("K19").PasteSpecial(-4122) | Out-Null

# Match the cursor position left after the data entry --------------------
$ws.Range("A25").Select() | Out-Null
